# v.0.7.0 74181 ALU added and working !

$wb = $excel.ActiveWorkbook

$wsControl = $wb.Worksheets.Item("Control Lines")
$wsTest = $wb.Worksheets.Item("Test Program")

# New cells in columns C:E use the sheet's "Text" number format (same as
# the existing rows in this table).
$wsTest.Range("C7:E9").NumberFormat = "@"

# Row 7 - ADD A, A, B
$wsTest.Range("A7").Value = "ADD A, A, B"
$wsTest.Range("F7").Value = "Sums A and B, put result in A"
$wsTest.Range("C7").Value = "14"
$wsTest.Range("D7").Value = "1x"
$wsTest.Range("E7").Value = "xxxx"

# Row 8 - SUB A, A, B
$wsTest.Range("C8").Value = "18"
$wsTest.Range("A8").Value = "SUB A, A, B"
$wsTest.Range("F8").Value = "A - B, result in A"
$wsTest.Range("D8").Value = "1x"
$wsTest.Range("E8").Value = "xxxx"

# Note about ALU inputs/outputs, on Control Lines sheet
$wsControl.Range("X7").Value = "All ALU operations take a fix register (A) as ALU input A and R2 as ALU input B. Output to R1"

# Row 9 - NOT A
$wsTest.Range("C9").Value = "1c"
$wsTest.Range("A9").Value = "NOT A"
$wsTest.Range("F9").Value = "NOT A, result in A"
$wsTest.Range("D9").Value = "0x"
$wsTest.Range("E9").Value = "xxxx"

# Rename "Test Program" sheet to "Some instructions for test"
$wsTest.Name = "Some instructions for test"

# Update selections to mirror authored state
$wsControl.Range("X8").Select() | Out-Null
$wsTest.Range("D9").Select() | Out-Null
$wsControl.Activate() | Out-Null
